$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.017.11"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.66%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.871.61"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.06%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4368"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -5.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3741"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.42%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07489"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.34%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9382"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.61%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.34"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.88%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.941.32"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.82%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.737"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.54%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.447"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.40%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06876"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.91%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.002"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.11%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "81.38"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.58%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009047"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.07%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.002"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.14%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.92"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.979.05"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.83%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.122"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.30%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.05"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.56%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.131.79"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.44%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.032"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.82%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "154.41"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.54"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.84%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.606"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.84%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "113.83"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.46%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.710"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -8.10%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09042"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.17%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.8160"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.823"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.91%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.180"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.20%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.977"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.22%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.001"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.13%  "
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.126"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.72%  "
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05524"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.50%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01978"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.61%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.963"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.90%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5275"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.57%  "
$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1704"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.95%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.022"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -7.08%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.797"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.06751"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.43%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4898"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.90%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.62"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.60%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "106.86"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.53%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.678"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.98%  "
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.914"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -14.21%  "
$ws.Range("B51").Value = "PaxDollar"
$ws.Range("C51").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.001"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.05%  "
